# Natmi following Dr Hou advice
# Rewrite the LR-pair data block (rows 2-9) on Sheet1 with the updated
# cell-type labels (ECs/M1/M2/sCs) and recomputed specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# One array per data row (row number, then values for columns A..T)
$dataRows = @(
    @(2, "ECs", "Ccl4", "Ccr3", "M1", 2, 0.6666666666666666, 0.08097766666666667, 0.242933, 0.0005396628443312377, 0.0005403049696870037, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 0.01072036336444444, 0.09648327027999999, 0.0001606706026779829, 0.0001608617788335946),
    @(3, "ECs", "Ccl4", "Ccr3", "M2", 2, 0.6666666666666666, 0.08097766666666667, 0.242933, 0.0005396628443312377, 0.0005403049696870037, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 0.02528735484344445, 0.227586193591, 0.0003789922416532548, 0.000379443190853409),
    @(4, "M1", "Ccl4", "Ccr3", "M1", 3, 1, 71.97958300000001, 215.938749, 0.4796965397318159, 0.4802673133443974, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 9.529137061426667, 85.76223355283999, 0.1428171921614588, 0.1429871251877723),
    @(5, "M1", "Ccl4", "Ccr3", "M2", 3, 1, 71.97958300000001, 215.938749, 0.4796965397318159, 0.4802673133443974, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 22.47747226771367, 202.297250409423, 0.3368793475703571, 0.337280188156625),
    @(6, "M2", "Ccl4", "Ccr3", "M1", 3, 1, 77.45677933333333, 232.370338, 0.5161984478983551, 0.5168126538150387, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 10.25424482667555, 92.28820344008, 0.15368468775731, 0.1538675516247016),
    @(7, "M2", "Ccl4", "Ccr3", "M2", 3, 1, 77.45677933333333, 232.370338, 0.5161984478983551, 0.5168126538150387, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 24.18786740416956, 217.690806637526, 0.3625137601410451, 0.3629451021903371),
    @(8, "sCs", "Ccl4", "Ccr3", "M1", 1, 0.5, 0.534989, 1.069978, 0.003565349525497869, 0.002379727870876994, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 0.07082541041333333, 0.42495246248, 0.001061490267556381, 0.0007085021976957098),
    @(9, "sCs", "Ccl4", "Ccr3", "M2", 1, 0.5, 0.534989, 1.069978, 0.003565349525497869, 0.002379727870876994, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 0.1670640466343334, 1.002384279806, 0.002503859257941488, 0.001671225673181284)
)

foreach ($dataRow in $dataRows) {
    $r = $dataRow[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $dataRow[$i + 1]
    }
}
